# Auto-generated edit script applying numeric corrections to Yojimbo_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 196.84
$ws.Range("I15").Value = 196.84
$ws.Range("K15").Value = 590.52
$ws.Range("M15").Value = -421.52

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2626
$ws.Range("I125").Value = 3800
$ws.Range("J125").Value = 2391.2
$ws.Range("K125").Value = 34200
$ws.Range("L125").Value = 21520.8
$ws.Range("M125").Value = -31740
$ws.Range("N125").Value = -26440.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 858.25
$ws.Range("J129").Value = 922.6667
$ws.Range("L129").Value = 2768.0001
$ws.Range("N129").Value = -12768.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2311.4883
$ws.Range("I141").Value = 2219.8438
$ws.Range("J141").Value = 2578.0908
$ws.Range("K141").Value = 6659.5314
$ws.Range("L141").Value = 7734.2724
$ws.Range("M141").Value = -1479.5314
$ws.Range("N141").Value = -18094.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2639.94
$ws.Range("I32").Value = 2387.0322
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 2387.0322
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -2100.0322
$ws.Range("N32").Value = -6574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3824
$ws.Range("I122").Value = 3102
$ws.Range("J122").Value = 5990
$ws.Range("K122").Value = 9306
$ws.Range("L122").Value = 17970
$ws.Range("M122").Value = -6856
$ws.Range("N122").Value = -22870

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2645.1455
$ws.Range("I132").Value = 2205.468
$ws.Range("J132").Value = 5228.25
$ws.Range("K132").Value = 6616.404
$ws.Range("L132").Value = 15684.75
$ws.Range("M132").Value = -4086.404
$ws.Range("N132").Value = -20744.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4149.362
$ws.Range("I86").Value = 3854.0977
$ws.Range("J86").Value = 6167
$ws.Range("K86").Value = 3854.0977
$ws.Range("L86").Value = 6167
$ws.Range("M86").Value = -2731.0977
$ws.Range("N86").Value = -8413

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4149.362
$ws.Range("I89").Value = 3854.0977
$ws.Range("J89").Value = 6167
$ws.Range("K89").Value = 19270.4885
$ws.Range("L89").Value = 30835
$ws.Range("M89").Value = -13654.4885
$ws.Range("N89").Value = -42067

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1134.8235
$ws.Range("I107").Value = 1064.7693
$ws.Range("J107").Value = 1362.5
$ws.Range("K107").Value = 1064.7693
$ws.Range("L107").Value = 1362.5
$ws.Range("M107").Value = 855.2307000000001
$ws.Range("N107").Value = -5202.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1154.8368
$ws.Range("I134").Value = 963.91895
$ws.Range("J134").Value = 1743.5
$ws.Range("K134").Value = 2891.75685
$ws.Range("L134").Value = 5230.5
$ws.Range("M134").Value = -356.7568499999998
$ws.Range("N134").Value = -10300.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 20056.572
$ws.Range("J50").Value = 20056.572
$ws.Range("L50").Value = 20056.572
$ws.Range("N50").Value = -21306.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1988.8214
$ws.Range("I58").Value = 2248.158
$ws.Range("J58").Value = 1441.3334
$ws.Range("K58").Value = 2248.158
$ws.Range("L58").Value = 1441.3334
$ws.Range("M58").Value = -2045.158
$ws.Range("N58").Value = -1847.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1687.4615
$ws.Range("I99").Value = 1270.7778
$ws.Range("J99").Value = 2625
$ws.Range("K99").Value = 1270.7778
$ws.Range("L99").Value = 2625
$ws.Range("M99").Value = 227.2221999999999
$ws.Range("N99").Value = -5621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 17686.46
$ws.Range("J108").Value = 17686.46
$ws.Range("L108").Value = 17686.46
$ws.Range("N108").Value = -25366.46

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1687.4615
$ws.Range("I126").Value = 1270.7778
$ws.Range("J126").Value = 2625
$ws.Range("K126").Value = 3812.3334
$ws.Range("L126").Value = 7875
$ws.Range("M126").Value = -1342.3334
$ws.Range("N126").Value = -12815

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1988.8214
$ws.Range("I136").Value = 2248.158
$ws.Range("J136").Value = 1441.3334
$ws.Range("K136").Value = 6744.474
$ws.Range("L136").Value = 4324.0002
$ws.Range("M136").Value = -4194.474
$ws.Range("N136").Value = -9424.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 22985.715
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 32140
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 32140
$ws.Range("M7").Value = 12
$ws.Range("N7").Value = -32364

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 22985.715
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 32140
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 32140
$ws.Range("M8").Value = 39
$ws.Range("N8").Value = -32418

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3216.4167
$ws.Range("I122").Value = 1901.75
$ws.Range("J122").Value = 3873.75
$ws.Range("K122").Value = 5705.25
$ws.Range("L122").Value = 11621.25
$ws.Range("M122").Value = -3255.25
$ws.Range("N122").Value = -16521.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1796.862
$ws.Range("I132").Value = 1494.4186
$ws.Range("J132").Value = 2663.8667
$ws.Range("K132").Value = 4483.2558
$ws.Range("L132").Value = 7991.6001
$ws.Range("M132").Value = -1953.2558
$ws.Range("N132").Value = -13051.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 504.8
$ws.Range("I16").Value = 504.8
$ws.Range("K16").Value = 504.8
$ws.Range("M16").Value = -334.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2675.8914
$ws.Range("I68").Value = 1538.4
$ws.Range("J68").Value = 2814.6099
$ws.Range("K68").Value = 1538.4
$ws.Range("L68").Value = 2814.6099
$ws.Range("M68").Value = -789.4000000000001
$ws.Range("N68").Value = -4312.609899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2675.8914
$ws.Range("I71").Value = 1538.4
$ws.Range("J71").Value = 2814.6099
$ws.Range("K71").Value = 7692
$ws.Range("L71").Value = 14073.0495
$ws.Range("M71").Value = -3948
$ws.Range("N71").Value = -21561.0495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2302.5
$ws.Range("I82").Value = 1353.5454
$ws.Range("J82").Value = 2720.04
$ws.Range("K82").Value = 1353.5454
$ws.Range("L82").Value = 2720.04
$ws.Range("M82").Value = -992.5454
$ws.Range("N82").Value = -3442.04

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2302.5
$ws.Range("I85").Value = 1353.5454
$ws.Range("J85").Value = 2720.04
$ws.Range("K85").Value = 1353.5454
$ws.Range("L85").Value = 2720.04
$ws.Range("M85").Value = -105.5454
$ws.Range("N85").Value = -5216.04

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1859.2941
$ws.Range("I93").Value = 1092.3334
$ws.Range("J93").Value = 3700
$ws.Range("K93").Value = 1092.3334
$ws.Range("L93").Value = 3700
$ws.Range("M93").Value = 155.6666
$ws.Range("N93").Value = -6196

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 779.2727
$ws.Range("I81").Value = 790
$ws.Range("J81").Value = 756.2857
$ws.Range("K81").Value = 1580
$ws.Range("L81").Value = 1512.5714
$ws.Range("M81").Value = -519
$ws.Range("N81").Value = -3634.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 779.2727
$ws.Range("I84").Value = 790
$ws.Range("J84").Value = 756.2857
$ws.Range("K84").Value = 7900
$ws.Range("L84").Value = 7562.857
$ws.Range("M84").Value = -2596
$ws.Range("N84").Value = -18170.857
